$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updated values
$ws.Range("B2").Value = 0.7287194209349384
$ws.Range("C2").Value = 86.29678392075563
$ws.Range("D2").Value = 3.082599426703578
$ws.Range("E2").Value = 71517.89157740913
$ws.Range("G2").Value = 71607.99968017753

# Row 3 updated values
$ws.Range("B3").Value = 0.02258322285507441
$ws.Range("C3").Value = 0.3375848360084654
$ws.Range("D3").Value = 116886.6739907443
$ws.Range("E3").Value = 71517.89157740913
$ws.Range("G3").Value = 188404.9257362123
